$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the static header value in A1 ("" -> "sn")
$ws.Range("A1").Value = "sn"

# Row 2: sn=1, A(old)=blank, B(old->C)=3, new D blank
$ws.Range("A2").Value = 1
$ws.Range("B2").ClearContents()
$ws.Range("B2").Font.Bold = $false
$ws.Range("C2").Value = 3
$ws.Range("D2").Font.Bold = $false

# Row 3: sn=2, 7, 8, 15
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 15

# Row 4: sn=3, 9, 9, blank
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 9
$ws.Range("D4").Font.Bold = $false
